$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-2.26%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'21"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'36.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.53%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'21"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.072"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.96%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'21"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.08059"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.79%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'21"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'1.976"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-2.54%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'21"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'7.803"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-2.18%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'21"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.9295"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.13%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'21"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.1464"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'35.39%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'21"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.1905"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.23%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'21"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.08965"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.82%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'21"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.03445"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.70%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'21"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.09860"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.59%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'21"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.001415"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.92%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'21"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.005775"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.17%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'21"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'3.532"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.67%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'21"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'4.052"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.84%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'21"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'2.834"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-3.90%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'21"
$ws.Range("G18").Style = "Normal"
$ws.Range("G19").Value = "'21"
$ws.Range("G19").Style = "Normal"
$ws.Range("E20").Value = "'-1.64%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'21"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'5.022"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.36%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'21"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.2392"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'9.28%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'21"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'0.04469"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.74%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'21"
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").Value = "'0.001205"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.63%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'21"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.004820"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.70%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'21"
$ws.Range("G25").Style = "Normal"
$ws.Range("E26").Value = "'-1.91%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'21"
$ws.Range("G26").Style = "Normal"
$ws.Range("E27").Value = "'-32.18%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'21"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'21"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'21"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'21"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'21"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'21"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'21"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'21"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'21"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'21"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'21"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'21"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.01895"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-4.57%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'21"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.04798"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-2.23%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'21"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.01058"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'8.28%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'21"
$ws.Range("G41").Style = "Normal"
$ws.Range("E42").Value = "'-6.29%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'21"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.1348"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.55%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'21"
$ws.Range("G43").Style = "Normal"
$ws.Range("E44").Value = "'-0.56%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'21"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.009726"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-15.84%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'21"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006212"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.42%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'21"
$ws.Range("G46").Style = "Normal"
$ws.Range("E47").Value = "'-0.30%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'21"
$ws.Range("G47").Style = "Normal"
$ws.Range("E48").Value = "'-0.09%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'21"
$ws.Range("G48").Style = "Normal"
$ws.Range("E49").Value = "'27.73%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'21"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002092"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.30%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'21"
$ws.Range("G50").Style = "Normal"
$ws.Range("E51").Value = "'-0.30%"
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").Value = "'21"
$ws.Range("G51").Style = "Normal"
